# Update NATMI derived values for rows 2-10 based on re-run analysis
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 149.656361
$ws.Range("H2").Value = 448.969083
$ws.Range("I2").Value = 0.5921360794347563
$ws.Range("J2").Value = 0.5921360794347564
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 17.16653
$ws.Range("N2").Value = 51.49959
$ws.Range("O2").Value = 0.0560345397128279
$ws.Range("P2").Value = 0.0560345397128279
$ws.Range("Q2").Value = 2569.080410797329
$ws.Range("R2").Value = 23121.72369717597
$ws.Range("S2").Value = 0.03318007265848506
$ws.Range("T2").Value = 0.03318007265848507

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 149.656361
$ws.Range("H3").Value = 448.969083
$ws.Range("I3").Value = 0.5921360794347563
$ws.Range("J3").Value = 0.5921360794347564
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 256.4443053333333
$ws.Range("N3").Value = 769.332916
$ws.Range("O3").Value = 0.8370788162388805
$ws.Range("P3").Value = 0.8370788162388805
$ws.Range("Q3").Value = 38378.52153535955
$ws.Range("R3").Value = 345406.6938182359
$ws.Range("S3").Value = 0.4956645684255775
$ws.Range("T3").Value = 0.4956645684255776

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 149.656361
$ws.Range("H4").Value = 448.969083
$ws.Range("I4").Value = 0.5921360794347563
$ws.Range("J4").Value = 0.5921360794347564
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 32.74538866666666
$ws.Range("N4").Value = 98.236166
$ws.Range("O4").Value = 0.1068866440482915
$ws.Range("P4").Value = 0.1068866440482915
$ws.Range("Q4").Value = 4900.555707383974
$ws.Range("R4").Value = 44105.00136645577
$ws.Range("S4").Value = 0.06329143835069367
$ws.Range("T4").Value = 0.06329143835069369

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 52.73412466666667
$ws.Range("H5").Value = 158.202374
$ws.Range("I5").Value = 0.208649853730866
$ws.Range("J5").Value = 0.208649853730866
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 17.16653
$ws.Range("N5").Value = 51.49959
$ws.Range("O5").Value = 0.0560345397128279
$ws.Range("P5").Value = 0.0560345397128279
$ws.Range("Q5").Value = 905.2619331140733
$ws.Range("R5").Value = 8147.357398026659
$ws.Range("S5").Value = 0.01169159851495794
$ws.Range("T5").Value = 0.01169159851495794

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 52.73412466666667
$ws.Range("H6").Value = 158.202374
$ws.Range("I6").Value = 0.208649853730866
$ws.Range("J6").Value = 0.208649853730866
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 256.4443053333333
$ws.Range("N6").Value = 769.332916
$ws.Range("O6").Value = 0.8370788162388805
$ws.Range("P6").Value = 0.8370788162388805
$ws.Range("Q6").Value = 13523.36596750473
$ws.Range("R6").Value = 121710.2937075426
$ws.Range("S6").Value = 0.1746563725694489
$ws.Range("T6").Value = 0.1746563725694489

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 52.73412466666667
$ws.Range("H7").Value = 158.202374
$ws.Range("I7").Value = 0.208649853730866
$ws.Range("J7").Value = 0.208649853730866
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 32.74538866666666
$ws.Range("N7").Value = 98.236166
$ws.Range("O7").Value = 0.1068866440482915
$ws.Range("P7").Value = 0.1068866440482915
$ws.Range("Q7").Value = 1726.799408206454
$ws.Range("R7").Value = 15541.19467385808
$ws.Range("S7").Value = 0.02230188264645917
$ws.Range("T7").Value = 0.02230188264645917

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 50.34932566666667
$ws.Range("H8").Value = 151.047977
$ws.Range("I8").Value = 0.1992140668343777
$ws.Range("J8").Value = 0.1992140668343777
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 17.16653
$ws.Range("N8").Value = 51.49959
$ws.Range("O8").Value = 0.0560345397128279
$ws.Range("P8").Value = 0.0560345397128279
$ws.Range("Q8").Value = 864.3232095366033
$ws.Range("R8").Value = 7778.90888582943
$ws.Range("S8").Value = 0.01116286853938489
$ws.Range("T8").Value = 0.01116286853938489

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 50.34932566666667
$ws.Range("H9").Value = 151.047977
$ws.Range("I9").Value = 0.1992140668343777
$ws.Range("J9").Value = 0.1992140668343777
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 256.4443053333333
$ws.Range("N9").Value = 769.332916
$ws.Range("O9").Value = 0.8370788162388805
$ws.Range("P9").Value = 0.8370788162388805
$ws.Range("Q9").Value = 12911.7978445901
$ws.Range("R9").Value = 116206.1806013109
$ws.Range("S9").Value = 0.1667578752438541
$ws.Range("T9").Value = 0.1667578752438541

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 50.34932566666667
$ws.Range("H10").Value = 151.047977
$ws.Range("I10").Value = 0.1992140668343777
$ws.Range("J10").Value = 0.1992140668343777
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 32.74538866666666
$ws.Range("N10").Value = 98.236166
$ws.Range("O10").Value = 0.1068866440482915
$ws.Range("P10").Value = 0.1068866440482915
$ws.Range("Q10").Value = 1648.708238059575
$ws.Range("R10").Value = 14838.37414253618
$ws.Range("S10").Value = 0.02129332305113868
$ws.Range("T10").Value = 0.02129332305113869
